$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.132.17'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.015.80'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.16'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.72'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.02%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0807'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.15'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.311.11'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.48%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.850'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.51'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.026.68'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.026.29'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.45'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.44'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.58'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.136'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.50%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.37'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.70%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0659'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.75%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.50'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.92%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0986'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.93'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0215'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.95'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.384.85'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.49'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.12'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +14.37%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.85'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.92'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.33%  '
